$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append (user_id, ingredient, added_at) starting at row 10.
# Copy the number-formatted style from an existing date cell (C9) down the
# new C-column cells first, so the new date values keep the same style
# index ("s=1", short date format) rather than Excel creating a brand new
# custom number format entry.
$ws.Range("C9").Copy($ws.Range("C10:C26"))

$rows = @(
    @(4, "eggs",       45996.82676960649),
    @(4, "chicken",    45996.82679936342),
    @(4, "fish",       45996.82682456019),
    @(4, "bacon",      45996.82685162037),
    @(4, "olive oil",  45996.8269062037),
    @(4, "pasta",      45996.82694400463),
    @(4, "tomato",     45996.82698706018),
    @(4, "potato",     45996.82705122685),
    @(4, "cheese",     45996.82707865741),
    @(4, "salt",       45996.82710900463),
    @(4, "pepper",     45996.82713510416),
    @(4, "tomatos",    45996.84235789352),
    @(5, "Peppers",    45996.85749776621),
    @(5, "chicken",    45996.8575313426),
    @(5, "egg",        45996.857547754626),
    @(5, "fish",       45996.857569872685),
    @(5, "potato",     45996.85764202546)
)

$r = 10
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
